$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.211.54'
$ws.Range("E2").Value = '  +0.82%  '

$ws.Range("D3").Value = '2.475.84'
$ws.Range("E3").Value = '  +0.72%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.33%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.14%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.508'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.55%  '

$ws.Range("D9").Value = '2.474.98'
$ws.Range("E9").Value = '  +0.59%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.152'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.02%  '

$ws.Range("E11").Value = '  +0.81%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.90'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.71%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.333'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.09%  '

$ws.Range("D14").Value = '69.101.06'
$ws.Range("E14").Value = '  +0.87%  '

$ws.Range("E15").Value = '  -0.25%  '

$ws.Range("E16").Value = '  -1.06%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '23.77'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.06%  '

$ws.Range("D18").Value = '2.418.86'
$ws.Range("E18").Value = '  -2.63%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.96%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '339.72'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.88%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.74%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.81'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.69%  '

$ws.Range("E24").Value = '  -0.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.34'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.54%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.70'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.68%  '

$ws.Range("D27").Value = '2.604.20'
$ws.Range("E27").Value = '  +0.98%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.31'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.82%  '

$ws.Range("D30").Value = '0.0₃0826'
$ws.Range("E30").Value = '  -1.92%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.22'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.70%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '433.71'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.06%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.15'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.23%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.63'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '157.68'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.07'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.25%  '

$ws.Range("B38").Value = 'USDe'
$ws.Range("C38").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.03%  '

$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.110'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.50%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.86'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.77%  '

$ws.Range("E41").Value = '  -1.54%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.46'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.38%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.48'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.67%  '

$ws.Range("E44").Value = '  -0.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.08'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.16%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '133.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.71%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.36'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.29%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0717'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.16%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.487'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.07%  '

$ws.Range("E50").Value = '  -0.17%  '

$ws.Range("E51").Value = '  +0.17%  '
